$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three shared-string values (B3, B5, B6) ---
# Force the cells to remain text (not auto-parsed as dates), then apply a
# cosmetic date-like custom number format (this mirrors the underlying
# OOXML change: new numFmt 165 "YYYY\-MM\-DD" applied as style index 1).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "april 2010"
$ws.Range("B3").NumberFormat = "YYYY\-MM\-DD"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "january 2011"
$ws.Range("B5").NumberFormat = "YYYY\-MM\-DD"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "april 2011"
$ws.Range("B6").NumberFormat = "General"

# --- Update the sheet selection ---
$ws.Range("B43").Select()
